$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1897.6666
$ws.Range("J43").Value = 1192.6666
$ws.Range("L43").Value = 1192.6666
$ws.Range("N43").Value = -1330.6666
$ws.Range("H64").Value = 4999.3335
$ws.Range("J64").Value = 4999
$ws.Range("L64").Value = 4999
$ws.Range("N64").Value = -5495
$ws.Range("H67").Value = 4999.3335
$ws.Range("J67").Value = 4999
$ws.Range("L67").Value = 4999
$ws.Range("N67").Value = -6715
$ws.Range("H100").Value = 3272
$ws.Range("I100").Value = 3546.4
$ws.Range("K100").Value = 3546.4
$ws.Range("M100").Value = -3005.4
$ws.Range("H112").Value = 4624.7
$ws.Range("I112").Value = 3375
$ws.Range("J112").Value = 4937.125
$ws.Range("K112").Value = 10125
$ws.Range("L112").Value = 14811.375
$ws.Range("M112").Value = -9017
$ws.Range("N112").Value = -17027.375
$ws.Range("H137").Value = 2696.0833
$ws.Range("J137").Value = 4449.25
$ws.Range("L137").Value = 13347.75
$ws.Range("N137").Value = -18447.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6212.2666
$ws.Range("I32").Value = 4483.4614
$ws.Range("J32").Value = 17449.5
$ws.Range("K32").Value = 4483.4614
$ws.Range("L32").Value = 17449.5
$ws.Range("M32").Value = -4196.4614
$ws.Range("N32").Value = -18023.5
$ws.Range("H61").Value = 3237.8462
$ws.Range("I61").Value = 2949.4
$ws.Range("K61").Value = 2949.4
$ws.Range("M61").Value = -2737.4
$ws.Range("H62").Value = 55000
$ws.Range("J62").Value = 55000
$ws.Range("L62").Value = 55000
$ws.Range("N62").Value = -56248
$ws.Range("H65").Value = 55000
$ws.Range("J65").Value = 55000
$ws.Range("L65").Value = 165000
$ws.Range("N65").Value = -171240
$ws.Range("H136").Value = 3237.8462
$ws.Range("I136").Value = 2949.4
$ws.Range("K136").Value = 8848.200000000001
$ws.Range("M136").Value = -6298.200000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 4073.8572
$ws.Range("I20").Value = 3433.3333
$ws.Range("J20").Value = 4554.25
$ws.Range("K20").Value = 3433.3333
$ws.Range("L20").Value = 4554.25
$ws.Range("M20").Value = -3186.3333
$ws.Range("N20").Value = -5048.25
$ws.Range("H134").Value = 4832.3335
$ws.Range("I134").Value = 4803.6113
$ws.Range("K134").Value = 14410.8339
$ws.Range("M134").Value = -11875.8339

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 1000
$ws.Range("I2").Value = 1000
$ws.Range("K2").Value = 1000
$ws.Range("M2").Value = -887
$ws.Range("H14").Value = 7997.5
$ws.Range("J14").Value = 7997.5
$ws.Range("L14").Value = 7997.5
$ws.Range("N14").Value = -8337.5
$ws.Range("H22").Value = 6667162
$ws.Range("J22").Value = 13333667
$ws.Range("L22").Value = 13333667
$ws.Range("N22").Value = -13334367
$ws.Range("H32").Value = 4350
$ws.Range("I32").Value = 4350
$ws.Range("K32").Value = 4350
$ws.Range("M32").Value = -4034
$ws.Range("H33").Value = 1705.6
$ws.Range("I33").Value = 1705.6
$ws.Range("K33").Value = 1705.6
$ws.Range("M33").Value = -1326.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H15").Value = 270.85715
$ws.Range("J15").Value = 98
$ws.Range("L15").Value = 294
$ws.Range("N15").Value = -574
$ws.Range("H16").Value = 530
$ws.Range("J16").Value = 530
$ws.Range("L16").Value = 1590
$ws.Range("N16").Value = -1936
$ws.Range("H70").Value = 12490.75
$ws.Range("I70").Value = 9987
$ws.Range("K70").Value = 29961
$ws.Range("M70").Value = -29646
$ws.Range("H73").Value = 12490.75
$ws.Range("I73").Value = 9987
$ws.Range("K73").Value = 29961
$ws.Range("M73").Value = -28869
$ws.Range("H122").Value = 560.7143
$ws.Range("I122").Value = 570.8333
$ws.Range("K122").Value = 5137.4997
$ws.Range("M122").Value = -2687.4997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H101").Value = 19999
$ws.Range("J101").Value = 19999
$ws.Range("L101").Value = 19999
$ws.Range("N101").Value = -26489
$ws.Range("H132").Value = 4670.4287
$ws.Range("I132").Value = 4364.4165
$ws.Range("K132").Value = 13093.2495
$ws.Range("M132").Value = -10563.2495

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 9082.429
$ws.Range("I7").Value = 9113
$ws.Range("J7").Value = 8899
$ws.Range("K7").Value = 9113
$ws.Range("L7").Value = 8899
$ws.Range("M7").Value = -9001
$ws.Range("N7").Value = -9123
$ws.Range("H40").Value = 3482.3333
$ws.Range("I40").Value = 3482.3333
$ws.Range("K40").Value = 3482.3333
$ws.Range("M40").Value = -3346.3333
$ws.Range("H46").Value = 1695.8334
$ws.Range("I46").Value = 507.5
$ws.Range("J46").Value = 4072.5
$ws.Range("K46").Value = 507.5
$ws.Range("L46").Value = 4072.5
$ws.Range("M46").Value = -319.5
$ws.Range("N46").Value = -4448.5
$ws.Range("H61").Value = 1711.4117
$ws.Range("I61").Value = 1749.9375
$ws.Range("J61").Value = 1095
$ws.Range("K61").Value = 1749.9375
$ws.Range("L61").Value = 1095
$ws.Range("M61").Value = -1547.9375
$ws.Range("N61").Value = -1499
$ws.Range("H82").Value = 2554
$ws.Range("I82").Value = 1766.6666
$ws.Range("K82").Value = 1766.6666
$ws.Range("M82").Value = -1405.6666
$ws.Range("H85").Value = 2554
$ws.Range("I85").Value = 1766.6666
$ws.Range("K85").Value = 1766.6666
$ws.Range("M85").Value = -518.6666
$ws.Range("H113").Value = 1711.4117
$ws.Range("I113").Value = 1749.9375
$ws.Range("J113").Value = 1095
$ws.Range("K113").Value = 1749.9375
$ws.Range("L113").Value = 1095
$ws.Range("M113").Value = 420.0625
$ws.Range("N113").Value = -5435
$ws.Range("H126").Value = 9082.429
$ws.Range("I126").Value = 9113
$ws.Range("J126").Value = 8899
$ws.Range("K126").Value = 27339
$ws.Range("L126").Value = 26697
$ws.Range("M126").Value = -24869
$ws.Range("N126").Value = -31637

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 4852.0713
$ws.Range("I122").Value = 5581.222
$ws.Range("J122").Value = 3539.6
$ws.Range("K122").Value = 16743.666
$ws.Range("L122").Value = 10618.8
$ws.Range("M122").Value = -14293.666
$ws.Range("N122").Value = -15518.8
